$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.077.46"
$ws.Range("E2").Value = "  -0.20%  "

$ws.Range("D3").Value = "1.655.62"
$ws.Range("E3").Value = "  +0.25%  "

$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.42%  "

$ws.Range("D5").Value = "218.48"
$ws.Range("E5").Value = "  +0.12%  "

$ws.Range("D6").Value = "0.5308"
$ws.Range("E6").Value = "  +2.08%  "

$ws.Range("E7").Value = "  -0.33%  "

$ws.Range("D8").Value = "0.2626"
$ws.Range("E8").Value = "  -1.31%  "

$ws.Range("D9").Value = "0.06320"
$ws.Range("E9").Value = "  +0.04%  "

$ws.Range("D10").Value = "20.47"
$ws.Range("E10").Value = "  -2.71%  "

$ws.Range("D11").Value = "0.07753"
$ws.Range("E11").Value = "  +0.30%  "

$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "4.497"
$ws.Range("E12").Value = "  +1.43%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.663.40"
$ws.Range("E13").Value = "  +0.82%  "

$ws.Range("D14").Value = "0.5482"
$ws.Range("E14").Value = "  +0.45%  "

$ws.Range("D15").Value = "0.0₅8117"
$ws.Range("E15").Value = "  -1.36%  "

$ws.Range("D16").Value = "65.23"
$ws.Range("E16").Value = "  +0.70%  "

$ws.Range("D17").Value = "26.093.48"
$ws.Range("E17").Value = "  -0.26%  "

$ws.Range("E18").Value = "  -0.34%  "

$ws.Range("D19").Value = "4.578"
$ws.Range("E19").Value = "  -1.92%  "

$ws.Range("D20").Value = "193.79"
$ws.Range("E20").Value = "  +0.69%  "

$ws.Range("D21").Value = "10.07"
$ws.Range("E21").Value = "  -0.76%  "

$ws.Range("D22").Value = "6.013"
$ws.Range("E22").Value = "  -1.27%  "

$ws.Range("E23").Value = "  -0.47%  "

$ws.Range("D24").Value = "139.87"
$ws.Range("E24").Value = "  +1.93%  "

$ws.Range("D25").Value = "0.1251"
$ws.Range("E25").Value = "  +0.93%  "

$ws.Range("D26").Value = "7.284"
$ws.Range("E26").Value = "  +0.81%  "

$ws.Range("D27").Value = "16.30"
$ws.Range("E27").Value = "  +1.16%  "

$ws.Range("D28").Value = "1.414"
$ws.Range("E28").Value = "  -1.03%  "

$ws.Range("D29").Value = "0.05966"
$ws.Range("E29").Value = "  -1.04%  "

$ws.Range("D30").Value = "1.278"
$ws.Range("E30").Value = "  -0.30%  "

$ws.Range("D31").Value = "3.510"
$ws.Range("E31").Value = "  -1.36%  "

$ws.Range("D32").Value = "3.257"
$ws.Range("E32").Value = "  -1.96%  "

$ws.Range("D33").Value = "1.548"
$ws.Range("E33").Value = "  -6.09%  "

$ws.Range("D34").Value = "2.414"
$ws.Range("E34").Value = "  +0.12%  "

$ws.Range("D35").Value = "0.9454"
$ws.Range("E35").Value = "  -3.38%  "

$ws.Range("D36").Value = "2.755"
$ws.Range("E36").Value = "  -0.53%  "

$ws.Range("D37").Value = "0.5651"
$ws.Range("E37").Value = "  -4.71%  "

$ws.Range("D38").Value = "0.01612"

$ws.Range("D39").Value = "5.886"
$ws.Range("E39").Value = "  -1.08%  "

$ws.Range("D40").Value = "0.8458"
$ws.Range("E40").Value = "  -2.02%  "

$ws.Range("E41").Value = "  -0.25%  "

$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "1.009.56"
$ws.Range("E42").Value = "  -2.93%  "

$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "100.91"
$ws.Range("E43").Value = "  +1.29%  "

$ws.Range("D44").Value = "1.792.10"
$ws.Range("E44").Value = "  -0.02%  "

$ws.Range("D45").Value = "57.11"
$ws.Range("E45").Value = "  +0.07%  "

$ws.Range("D46").Value = "0.0₈106"
$ws.Range("E46").Value = "  -4.09%  "

$ws.Range("D47").Value = "1.007"
$ws.Range("E47").Value = "  +0.16%  "

$ws.Range("D48").Value = "1.489"
$ws.Range("E48").Value = "  +1.44%  "

$ws.Range("D49").Value = "0.4290"
$ws.Range("E49").Value = "  +1.44%  "

$ws.Range("D50").Value = "7.845"
$ws.Range("E50").Value = "  -3.17%  "
